$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(31, 8).Value = 470
$ws.Cells.Item(31, 9).Value = 470
$ws.Cells.Item(31, 11).Value = 1410
$ws.Cells.Item(31, 13).Value = -1180
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3926.3333
$ws.Cells.Item(32, 9).Value = 3569.0588
$ws.Cells.Item(32, 11).Value = 3569.0588
$ws.Cells.Item(32, 13).Value = -3282.0588
$ws.Cells.Item(45, 8).Value = 9555.223
$ws.Cells.Item(45, 10).Value = 8600
$ws.Cells.Item(45, 12).Value = 8600
$ws.Cells.Item(45, 14).Value = -9354
$ws.Cells.Item(76, 8).Value = 70000
$ws.Cells.Item(76, 10).Value = 70000
$ws.Cells.Item(76, 12).Value = 70000
$ws.Cells.Item(76, 14).Value = -70676
$ws.Cells.Item(79, 8).Value = 70000
$ws.Cells.Item(79, 10).Value = 70000
$ws.Cells.Item(79, 12).Value = 70000
$ws.Cells.Item(79, 14).Value = -72340
$ws.Cells.Item(80, 8).Value = 75749.164
$ws.Cells.Item(80, 10).Value = 84899
$ws.Cells.Item(80, 12).Value = 84899
$ws.Cells.Item(80, 14).Value = -86895
$ws.Cells.Item(83, 8).Value = 75749.164
$ws.Cells.Item(83, 10).Value = 84899
$ws.Cells.Item(83, 12).Value = 254697
$ws.Cells.Item(83, 14).Value = -264681
$ws.Cells.Item(110, 8).Value = 1661.931
$ws.Cells.Item(110, 9).Value = 850.5909
$ws.Cells.Item(110, 10).Value = 4211.857
$ws.Cells.Item(110, 11).Value = 850.5909
$ws.Cells.Item(110, 12).Value = 4211.857
$ws.Cells.Item(110, 13).Value = 1194.4091
$ws.Cells.Item(110, 14).Value = -8301.857
$ws.Cells.Item(132, 8).Value = 3229.2222
$ws.Cells.Item(132, 9).Value = 1747.0526
$ws.Cells.Item(132, 11).Value = 5241.1578
$ws.Cells.Item(132, 13).Value = -2711.1578
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 94872.5
$ws.Cells.Item(140, 10).Value = 94872.5
$ws.Cells.Item(140, 12).Value = 94872.5
$ws.Cells.Item(140, 14).Value = -105232.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 764.8889
$ws.Cells.Item(22, 9).Value = 447.25
$ws.Cells.Item(22, 10).Value = 1019
$ws.Cells.Item(22, 11).Value = 447.25
$ws.Cells.Item(22, 12).Value = 1019
$ws.Cells.Item(22, 13).Value = -97.25
$ws.Cells.Item(22, 14).Value = -1719
$ws.Cells.Item(58, 8).Value = 3097.3125
$ws.Cells.Item(58, 9).Value = 3013.087
$ws.Cells.Item(58, 10).Value = 3312.5557
$ws.Cells.Item(58, 11).Value = 3013.087
$ws.Cells.Item(58, 12).Value = 3312.5557
$ws.Cells.Item(58, 13).Value = -2810.087
$ws.Cells.Item(58, 14).Value = -3718.5557
$ws.Cells.Item(68, 8).Value = 54998
$ws.Cells.Item(68, 10).Value = 54998
$ws.Cells.Item(68, 12).Value = 54998
$ws.Cells.Item(68, 14).Value = -56496
$ws.Cells.Item(70, 8).Value = 47313
$ws.Cells.Item(70, 10).Value = 47313
$ws.Cells.Item(70, 12).Value = 47313
$ws.Cells.Item(70, 14).Value = -47943
$ws.Cells.Item(71, 8).Value = 54998
$ws.Cells.Item(71, 10).Value = 54998
$ws.Cells.Item(71, 12).Value = 164994
$ws.Cells.Item(71, 14).Value = -172482
$ws.Cells.Item(73, 8).Value = 47313
$ws.Cells.Item(73, 10).Value = 47313
$ws.Cells.Item(73, 12).Value = 47313
$ws.Cells.Item(73, 14).Value = -49497
$ws.Cells.Item(122, 8).Value = 25759.2
$ws.Cells.Item(122, 10).Value = 2197
$ws.Cells.Item(122, 12).Value = 6591
$ws.Cells.Item(122, 14).Value = -11491
$ws.Cells.Item(134, 8).Value = 1588.5555
$ws.Cells.Item(134, 9).Value = 1439.6666
$ws.Cells.Item(134, 11).Value = 4318.9998
$ws.Cells.Item(134, 13).Value = -1783.9998
$ws.Cells.Item(136, 8).Value = 3097.3125
$ws.Cells.Item(136, 9).Value = 3013.087
$ws.Cells.Item(136, 10).Value = 3312.5557
$ws.Cells.Item(136, 11).Value = 9039.261
$ws.Cells.Item(136, 12).Value = 9937.667099999999
$ws.Cells.Item(136, 13).Value = -6489.261
$ws.Cells.Item(136, 14).Value = -15037.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 668528.8
$ws.Cells.Item(5, 9).Value = 682.625
$ws.Cells.Item(5, 10).Value = 1431781.6
$ws.Cells.Item(5, 11).Value = 2047.875
$ws.Cells.Item(5, 12).Value = 4295344.800000001
$ws.Cells.Item(5, 13).Value = -1935.875
$ws.Cells.Item(5, 14).Value = -4295568.800000001
$ws.Cells.Item(33, 8).Value = 38.214287
$ws.Cells.Item(33, 10).Value = 54
$ws.Cells.Item(33, 12).Value = 324
$ws.Cells.Item(33, 14).Value = -890
$ws.Cells.Item(68, 8).Value = 11999.583
$ws.Cells.Item(68, 10).Value = 15393.889
$ws.Cells.Item(68, 12).Value = 46181.667
$ws.Cells.Item(68, 14).Value = -47803.667
$ws.Cells.Item(71, 8).Value = 11999.583
$ws.Cells.Item(71, 10).Value = 15393.889
$ws.Cells.Item(71, 12).Value = 138545.001
$ws.Cells.Item(71, 14).Value = -146657.001
$ws.Cells.Item(80, 8).Value = 68013.53999999999
$ws.Cells.Item(80, 9).Value = 4150
$ws.Cells.Item(80, 10).Value = 79625.09
$ws.Cells.Item(80, 11).Value = 12450
$ws.Cells.Item(80, 12).Value = 238875.27
$ws.Cells.Item(80, 13).Value = -11514
$ws.Cells.Item(80, 14).Value = -240747.27
$ws.Cells.Item(83, 8).Value = 68013.53999999999
$ws.Cells.Item(83, 9).Value = 4150
$ws.Cells.Item(83, 10).Value = 79625.09
$ws.Cells.Item(83, 11).Value = 37350
$ws.Cells.Item(83, 12).Value = 716625.8099999999
$ws.Cells.Item(83, 13).Value = -32670
$ws.Cells.Item(83, 14).Value = -725985.8099999999
$ws.Cells.Item(92, 8).Value = 579.6
$ws.Cells.Item(92, 9).Value = 600
$ws.Cells.Item(92, 10).Value = 566
$ws.Cells.Item(92, 11).Value = 1800
$ws.Cells.Item(92, 12).Value = 1698
$ws.Cells.Item(92, 13).Value = -552
$ws.Cells.Item(92, 14).Value = -4194
$ws.Cells.Item(135, 8).Value = 668528.8
$ws.Cells.Item(135, 9).Value = 682.625
$ws.Cells.Item(135, 10).Value = 1431781.6
$ws.Cells.Item(135, 11).Value = 6143.625
$ws.Cells.Item(135, 12).Value = 12886034.4
$ws.Cells.Item(135, 13).Value = -3608.625
$ws.Cells.Item(135, 14).Value = -12891104.4
$ws.Cells.Item(137, 8).Value = 7043.3125
$ws.Cells.Item(137, 10).Value = 13463.714
$ws.Cells.Item(137, 12).Value = 40391.142
$ws.Cells.Item(137, 14).Value = -50591.142
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 8622.825999999999
$ws.Cells.Item(102, 9).Value = 9295.388999999999
$ws.Cells.Item(102, 11).Value = 9295.388999999999
$ws.Cells.Item(102, 13).Value = -7673.388999999999
$ws.Cells.Item(132, 8).Value = 3839.5
$ws.Cells.Item(132, 9).Value = 3987.5757
$ws.Cells.Item(132, 11).Value = 11962.7271
$ws.Cells.Item(132, 13).Value = -9432.7271
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(63, 8).Value = 59173.332
$ws.Cells.Item(66, 8).Value = 59173.332
$ws.Cells.Item(100, 8).Value = 10389.8
$ws.Cells.Item(100, 9).Value = 10235.571
$ws.Cells.Item(100, 10).Value = 10749.667
$ws.Cells.Item(100, 11).Value = 10235.571
$ws.Cells.Item(100, 12).Value = 10749.667
$ws.Cells.Item(100, 13).Value = -9694.571
$ws.Cells.Item(100, 14).Value = -11831.667
$ws.Cells.Item(122, 8).Value = 6788.913
$ws.Cells.Item(122, 9).Value = 7708.8887
$ws.Cells.Item(122, 11).Value = 23126.6661
$ws.Cells.Item(122, 13).Value = -20676.6661
$ws.Cells.Item(132, 8).Value = 517595.8
$ws.Cells.Item(132, 9).Value = 1067120.4
$ws.Cells.Item(132, 10).Value = 4706.2
$ws.Cells.Item(132, 11).Value = 3201361.2
$ws.Cells.Item(132, 12).Value = 14118.6
$ws.Cells.Item(132, 13).Value = -3198831.2
$ws.Cells.Item(132, 14).Value = -19178.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 233194.6
$ws.Cells.Item(62, 9).Value = 428981.75
$ws.Cells.Item(62, 11).Value = 428981.75
$ws.Cells.Item(62, 13).Value = -428357.75
$ws.Cells.Item(65, 8).Value = 233194.6
$ws.Cells.Item(65, 9).Value = 428981.75
$ws.Cells.Item(65, 11).Value = 2144908.75
$ws.Cells.Item(65, 13).Value = -2141788.75
$ws.Cells.Item(70, 8).Value = 587972.5
$ws.Cells.Item(70, 10).Value = 773966.7
$ws.Cells.Item(70, 12).Value = 773966.7
$ws.Cells.Item(70, 14).Value = -774596.7
$ws.Cells.Item(73, 8).Value = 587972.5
$ws.Cells.Item(73, 10).Value = 773966.7
$ws.Cells.Item(73, 12).Value = 773966.7
$ws.Cells.Item(73, 14).Value = -776150.7
$ws.Cells.Item(81, 8).Value = 9171.846
$ws.Cells.Item(81, 9).Value = 13455.5
$ws.Cells.Item(81, 10).Value = 2318
$ws.Cells.Item(81, 11).Value = 26911
$ws.Cells.Item(81, 12).Value = 4636
$ws.Cells.Item(81, 13).Value = -25850
$ws.Cells.Item(81, 14).Value = -6758
$ws.Cells.Item(84, 8).Value = 9171.846
$ws.Cells.Item(84, 9).Value = 13455.5
$ws.Cells.Item(84, 10).Value = 2318
$ws.Cells.Item(84, 11).Value = 134555
$ws.Cells.Item(84, 12).Value = 23180
$ws.Cells.Item(84, 13).Value = -129251
$ws.Cells.Item(84, 14).Value = -33788
$ws.Cells.Item(100, 8).Value = 43872.223
$ws.Cells.Item(100, 9).Value = 10692.857
$ws.Cells.Item(100, 11).Value = 21385.714
$ws.Cells.Item(100, 13).Value = -20844.714
$ws.Cells.Item(122, 8).Value = 5143.7837
$ws.Cells.Item(122, 9).Value = 3410.8823
$ws.Cells.Item(122, 11).Value = 10232.6469
$ws.Cells.Item(122, 13).Value = -7782.6469
$ws.Cells.Item(132, 8).Value = 5975.278
$ws.Cells.Item(132, 10).Value = 5434.6
$ws.Cells.Item(132, 12).Value = 16303.8
$ws.Cells.Item(132, 14).Value = -21363.8
